$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (rows 2-11) from 1 to 2
$ws.Range("A2:A11").Value = 2

# Update the selection to T7
$ws.Range("T7").Select()
